# Insert a new data row before the current row 63 (shifts old rows 63-104
# down to 64-105, growing the used range to A1:R105), then populate the
# new row 63 with the new "Haba" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(63).Insert()

$ws.Range("A63").Value = 5
$ws.Range("B63").Value = "Macroferia Regional de Talca"
$ws.Range("C63").Value = "Maule"
$ws.Range("D63").Value = 44873
$ws.Range("E63").Value = 7
$ws.Range("F63").Value = 100112026
$ws.Range("G63").Value = "Haba"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 500
$ws.Range("K63").Value = 7000
$ws.Range("L63").Value = 7000
$ws.Range("M63").Value = 7000
$ws.Range("N63").Value = "$/saco 25 kilos"
$ws.Range("O63").Value = "Región del Maule"
$ws.Range("P63").Value = 280
$ws.Range("Q63").Value = 25
$ws.Range("R63").Value = "Hortaliza"
